$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Hydrogen): corrected Iron & steel value, Non-metallic minerals value removed
$ws.Range("B3").Value = 3073889.306338502
$ws.Range("D3").Value = ""

# Row 4 (Methanol): corrected Chemicals value
$ws.Range("C4").Value = 27.34422306158211

# Row 5 (Ammonia): corrected Chemicals value
$ws.Range("C5").Value = 482.3069101150899

# Row 7: rename label "Other" -> "Biogas" and update its value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 2076.146877920774

# New row 8: "Other" row, carrying over the label's formatting from row 7
$ws.Range("A7").Copy($ws.Range("A8"))
$ws.Range("A8").Value = "Other"
$ws.Range("D8").Value = 724.6697314911828
